$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.687922358512878
$ws.Range("B1").Value = 1.78466010093689
$ws.Range("C1").Value = 1.940633296966553
$ws.Range("D1").Value = 2.81668496131897
$ws.Range("E1").Value = 3.852812767028809
